$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are stored as text in this sheet (matching the existing cells),
# so force each written cell to Text format before assigning the new
# numeric-looking value. This avoids Excel auto-converting "45" etc. to
# a numeric value, keeping the cell type consistent with the rest of the sheet.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "45"
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "10"
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "19"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "34"
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "26"
$ws.Range("F1").NumberFormat = "@"
$ws.Range("F1").Value = "38"
$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "21"
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").Value = "12"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "41"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "48"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "40"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "4"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "50"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "27"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "32"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "15"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "17"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "7"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "46"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "16"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "47"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "49"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "29"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "31"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "23"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "39"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "52"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "35"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "33"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "5"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "37"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "53"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "18"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "57"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "11"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "56"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "24"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "36"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "42"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "13"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "22"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "30"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "8"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "14"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "6"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "20"
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "3"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "44"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "25"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "43"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "55"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "51"
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "1"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "58"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "54"
